$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (2021-09-02 .. 2021-09-09), appended after the existing
# last row (366) of the sheet.
$newRows = @(
    @{ Row = 367; A = 44441; B = 0; C = 13; D = 216.2701713525204 },
    @{ Row = 368; A = 44442; B = 2; C = 8;  D = 133.0893362169356 },
    @{ Row = 369; A = 44443; B = 0; C = 8;  D = 133.0893362169356 },
    @{ Row = 370; A = 44444; B = 1; C = 8;  D = 133.0893362169356 },
    @{ Row = 371; A = 44445; B = 0; C = 6;  D = 99.81700216270171 },
    @{ Row = 372; A = 44446; B = 1; C = 4;  D = 66.54466810846782 },
    @{ Row = 373; A = 44447; B = 0; C = 4;  D = 66.54466810846782 },
    @{ Row = 374; A = 44448; B = 1; C = 5;  D = 83.18083513558476 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Column A carries the date style (s="2" in the source workbook) that
    # the preceding rows use. Copy just the formatting from the row above
    # so the shared style index is reused instead of minting a new one.
    $ws.Range("A" + ($rowNum - 1)).Copy()
    $ws.Range("A" + $rowNum).PasteSpecial(-4122)

    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
}

$excel.CutCopyMode = 0
